# Auto-generated edit script applying numeric updates per the commit diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 339.92856
$ws.Range("I12").Value = 345.44446
$ws.Range("K12").Value = 345.44446
$ws.Range("M12").Value = -175.44446
$ws.Range("H18").Value = 566.6667
$ws.Range("I18").Value = 566.6667
$ws.Range("K18").Value = 566.6667
$ws.Range("M18").Value = -282.6667
$ws.Range("H64").Value = 4131.4
$ws.Range("I64").Value = 3915.2942
$ws.Range("J64").Value = 5356
$ws.Range("K64").Value = 3915.2942
$ws.Range("L64").Value = 5356
$ws.Range("M64").Value = -3667.2942
$ws.Range("N64").Value = -5852
$ws.Range("H67").Value = 4131.4
$ws.Range("I67").Value = 3915.2942
$ws.Range("J67").Value = 5356
$ws.Range("K67").Value = 3915.2942
$ws.Range("L67").Value = 5356
$ws.Range("M67").Value = -3057.2942
$ws.Range("N67").Value = -7072
$ws.Range("H98").Value = 786.34375
$ws.Range("I98").Value = 786.34375
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 786.34375
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = 711.65625
$ws.Range("N98").ClearContents()
$ws.Range("H122").Value = 786.34375
$ws.Range("I122").Value = 786.34375
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 2359.03125
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = 90.96875
$ws.Range("N122").ClearContents()
$ws.Range("H131").Value = 2256.9
$ws.Range("I131").Value = 1538.4286
$ws.Range("J131").Value = 3933.3333
$ws.Range("K131").Value = 4615.2858
$ws.Range("L131").Value = 11799.9999
$ws.Range("M131").Value = 424.7142000000003
$ws.Range("N131").Value = -21879.9999
$ws.Range("H132").Value = 1197.6666
$ws.Range("I132").Value = 692.9583
$ws.Range("J132").Value = 5235.3335
$ws.Range("K132").Value = 2078.8749
$ws.Range("L132").Value = 15706.0005
$ws.Range("M132").Value = 451.1251000000002
$ws.Range("N132").Value = -20766.0005
$ws.Range("H138").Value = 2803.675
$ws.Range("I138").Value = 925.2381
$ws.Range("J138").Value = 4879.8423
$ws.Range("K138").Value = 2775.7143
$ws.Range("L138").Value = 14639.5269
$ws.Range("M138").Value = 2364.2857
$ws.Range("N138").Value = -24919.5269

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2441.3125
$ws.Range("I2").Value = 3070.1
$ws.Range("J2").Value = 1393.3334
$ws.Range("K2").Value = 3070.1
$ws.Range("L2").Value = 1393.3334
$ws.Range("M2").Value = -2957.1
$ws.Range("N2").Value = -1619.3334
$ws.Range("H32").Value = 4172.0674
$ws.Range("I32").Value = 2985.1147
$ws.Range("K32").Value = 2985.1147
$ws.Range("M32").Value = -2698.1147
$ws.Range("H61").Value = 4866.7812
$ws.Range("I61").Value = 5176.3213
$ws.Range("J61").Value = 2700
$ws.Range("K61").Value = 5176.3213
$ws.Range("L61").Value = 2700
$ws.Range("M61").Value = -4964.3213
$ws.Range("N61").Value = -3124
$ws.Range("H63").Value = 100002690
$ws.Range("I63").Value = 142859550
$ws.Range("J63").Value = 3333
$ws.Range("K63").Value = 142859550
$ws.Range("L63").Value = 3333
$ws.Range("M63").Value = -142858864
$ws.Range("N63").Value = -4705
$ws.Range("H66").Value = 100002690
$ws.Range("I66").Value = 142859550
$ws.Range("J66").Value = 3333
$ws.Range("K66").Value = 714297750
$ws.Range("L66").Value = 16665
$ws.Range("M66").Value = -714294318
$ws.Range("N66").Value = -23529
$ws.Range("H116").Value = 2441.3125
$ws.Range("I116").Value = 3070.1
$ws.Range("J116").Value = 1393.3334
$ws.Range("K116").Value = 3070.1
$ws.Range("L116").Value = 1393.3334
$ws.Range("M116").Value = -776.0999999999999
$ws.Range("N116").Value = -5981.3334
$ws.Range("H123").Value = 44424.332
$ws.Range("J123").Value = 44424.332
$ws.Range("L123").Value = 44424.332
$ws.Range("N123").Value = -54224.332
$ws.Range("H132").Value = 3528.1177
$ws.Range("I132").Value = 1740.0476
$ws.Range("J132").Value = 6416.5386
$ws.Range("K132").Value = 5220.142800000001
$ws.Range("L132").Value = 19249.6158
$ws.Range("M132").Value = -2690.142800000001
$ws.Range("N132").Value = -24309.6158
$ws.Range("H135").Value = 71332.25
$ws.Range("J135").Value = 71332.25
$ws.Range("L135").Value = 71332.25
$ws.Range("N135").Value = -81472.25
$ws.Range("H136").Value = 4866.7812
$ws.Range("I136").Value = 5176.3213
$ws.Range("J136").Value = 2700
$ws.Range("K136").Value = 15528.9639
$ws.Range("L136").Value = 8100
$ws.Range("M136").Value = -12978.9639
$ws.Range("N136").Value = -13200

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2441.3125
$ws.Range("I3").Value = 3070.1
$ws.Range("J3").Value = 1393.3334
$ws.Range("K3").Value = 3070.1
$ws.Range("L3").Value = 1393.3334
$ws.Range("M3").Value = -2956.1
$ws.Range("N3").Value = -1621.3334
$ws.Range("H132").Value = 38000
$ws.Range("J132").Value = 38000
$ws.Range("L132").Value = 38000
$ws.Range("N132").Value = -48120
$ws.Range("H134").Value = 5857.2256
$ws.Range("I134").Value = 7656.7896
$ws.Range("J134").Value = 3007.9167
$ws.Range("K134").Value = 22970.3688
$ws.Range("L134").Value = 9023.750100000001
$ws.Range("M134").Value = -20435.3688
$ws.Range("N134").Value = -14093.7501

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H17").Value = 166708340
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 166708340
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 166708340
$ws.Range("M17").ClearContents()
$ws.Range("N17").Value = -166708688
$ws.Range("H41").Value = 9014.75
$ws.Range("I41").Value = 1529.5
$ws.Range("J41").Value = 16500
$ws.Range("K41").Value = 1529.5
$ws.Range("L41").Value = 16500
$ws.Range("M41").Value = -1101.5
$ws.Range("N41").Value = -17356
$ws.Range("H50").Value = 19996.666
$ws.Range("J50").Value = 19996.666
$ws.Range("L50").Value = 19996.666
$ws.Range("N50").Value = -21246.666
$ws.Range("H51").Value = 17481.666
$ws.Range("I51").Value = 10000
$ws.Range("J51").Value = 18978
$ws.Range("K51").Value = 10000
$ws.Range("L51").Value = 18978
$ws.Range("M51").Value = -9264
$ws.Range("N51").Value = -20450
$ws.Range("H59").Value = 25719.6
$ws.Range("I59").Value = 15000
$ws.Range("J59").Value = 28399.5
$ws.Range("K59").Value = 15000
$ws.Range("L59").Value = 28399.5
$ws.Range("M59").Value = -13855
$ws.Range("N59").Value = -30689.5
$ws.Range("H60").Value = 12054.074
$ws.Range("J60").Value = 12054.074
$ws.Range("L60").Value = 12054.074
$ws.Range("N60").Value = -13076.074
$ws.Range("H61").Value = 17481.666
$ws.Range("I61").Value = 10000
$ws.Range("J61").Value = 18978
$ws.Range("K61").Value = 10000
$ws.Range("L61").Value = 18978
$ws.Range("M61").Value = -9652
$ws.Range("N61").Value = -19674
$ws.Range("H62").Value = 8000
$ws.Range("I62").Value = 20000
$ws.Range("K62").Value = 20000
$ws.Range("M62").Value = -19376
$ws.Range("H65").Value = 8000
$ws.Range("I65").Value = 20000
$ws.Range("K65").Value = 100000
$ws.Range("M65").Value = -96880
$ws.Range("H68").Value = 30149.5
$ws.Range("J68").Value = 30149.5
$ws.Range("L68").Value = 30149.5
$ws.Range("N68").Value = -31647.5
$ws.Range("H71").Value = 30149.5
$ws.Range("J71").Value = 30149.5
$ws.Range("L71").Value = 90448.5
$ws.Range("N71").Value = -97936.5
$ws.Range("H74").Value = 21175.2
$ws.Range("J74").Value = 21175.2
$ws.Range("L74").Value = 21175.2
$ws.Range("N74").Value = -22923.2
$ws.Range("H77").Value = 21175.2
$ws.Range("J77").Value = 21175.2
$ws.Range("L77").Value = 63525.60000000001
$ws.Range("N77").Value = -72261.60000000001
$ws.Range("H132").Value = 2987.6
$ws.Range("I132").Value = 2703.6
$ws.Range("J132").Value = 3839.6
$ws.Range("K132").Value = 8110.799999999999
$ws.Range("L132").Value = 11518.8
$ws.Range("M132").Value = -5580.799999999999
$ws.Range("N132").Value = -16578.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H70").Value = 3911.111
$ws.Range("I70").Value = 2100
$ws.Range("J70").Value = 4428.5713
$ws.Range("K70").Value = 6300
$ws.Range("L70").Value = 13285.7139
$ws.Range("M70").Value = -5985
$ws.Range("N70").Value = -13915.7139
$ws.Range("H73").Value = 3911.111
$ws.Range("I73").Value = 2100
$ws.Range("J73").Value = 4428.5713
$ws.Range("K73").Value = 6300
$ws.Range("L73").Value = 13285.7139
$ws.Range("M73").Value = -5208
$ws.Range("N73").Value = -15469.7139
$ws.Range("H76").Value = 2854.3333
$ws.Range("I76").Value = 1708.6666
$ws.Range("J76").Value = 4000
$ws.Range("K76").Value = 5125.9998
$ws.Range("L76").Value = 12000
$ws.Range("M76").Value = -4742.9998
$ws.Range("N76").Value = -12766
$ws.Range("H79").Value = 2854.3333
$ws.Range("I79").Value = 1708.6666
$ws.Range("J79").Value = 4000
$ws.Range("K79").Value = 5125.9998
$ws.Range("L79").Value = 12000
$ws.Range("M79").Value = -3799.9998
$ws.Range("N79").Value = -14652
$ws.Range("H98").Value = 12500388
$ws.Range("I98").Value = 149
$ws.Range("J98").Value = 16667134
$ws.Range("K98").Value = 447
$ws.Range("L98").Value = 50001402
$ws.Range("M98").Value = 1051
$ws.Range("N98").Value = -50004398

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1993
$ws.Range("I102").Value = 2053.3333
$ws.Range("J102").Value = 1450
$ws.Range("K102").Value = 2053.3333
$ws.Range("L102").Value = 1450
$ws.Range("M102").Value = -431.3332999999998
$ws.Range("N102").Value = -4694

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 17544766
$ws.Range("I46").Value = 37037800
$ws.Range("K46").Value = 37037800
$ws.Range("M46").Value = -37037612
$ws.Range("H136").Value = 5822.8
$ws.Range("I136").Value = 6435.96
$ws.Range("J136").Value = 2757
$ws.Range("K136").Value = 19307.88
$ws.Range("L136").Value = 8271
$ws.Range("M136").Value = -16757.88
$ws.Range("N136").Value = -13371

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 50557.125
$ws.Range("J46").Value = 50557.125
$ws.Range("L46").Value = 50557.125
$ws.Range("N46").Value = -51019.125
$ws.Range("H81").Value = 1000
$ws.Range("J81").Value = 1000
$ws.Range("L81").Value = 2000
$ws.Range("N81").Value = -4122
$ws.Range("H84").Value = 1000
$ws.Range("J84").Value = 1000
$ws.Range("L84").Value = 10000
$ws.Range("N84").Value = -20608
$ws.Range("H134").Value = 50557.125
$ws.Range("J134").Value = 50557.125
$ws.Range("L134").Value = 151671.375
$ws.Range("N134").Value = -156741.375
